$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 6; $r++) {
    $ws.Cells.Item($r, 3).Value = "Error: unprocessable date"  # C
    $ws.Cells.Item($r, 4).Value = "Error: unprocessable date"  # D
    $ws.Cells.Item($r, 5).Value = "Error: not a number"        # E
    $ws.Cells.Item($r, 6).Value = "7"                          # F
}

$ws.Range("E2:E6").Select() | Out-Null
